$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LH_CRS")

$ws.Range("B2").Value = "LH-CRS-REGISTRATION-001"
$ws.Range("B3").Value = "LH-CRS-NAVIGATION-002"
$ws.Range("B4").Value = "LH-CRS-PUBLISH-003"
